# Updated cryptos list on Fri May 12 16:03:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.420.27'
$ws.Range('E2').Value = '''  -2.94%  '
$ws.Range('D3').Value = '''1.772.55'
$ws.Range('E3').Value = '''  -2.10%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E5').Value = '''  +0.21%  '
$ws.Range('D6').Value = '''306.60'
$ws.Range('E6').Value = '''  -1.16%  '
$ws.Range('D7').Value = '''0.4261'
$ws.Range('E7').Value = '''  +1.36%  '
$ws.Range('D8').Value = '''0.3609'
$ws.Range('E8').Value = '''  +1.55%  '
$ws.Range('D9').Value = '''0.07136'
$ws.Range('E9').Value = '''  +0.55%  '
$ws.Range('D10').Value = '''0.8382'
$ws.Range('E10').Value = '''  -1.38%  '
$ws.Range('D11').Value = '''20.38'
$ws.Range('E11').Value = '''  +0.87%  '
$ws.Range('D12').Value = '''1.798.20'
$ws.Range('E12').Value = '''  +0.18%  '
$ws.Range('B13').Value = '''Polkadot'
$ws.Range('C13').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.243'
$ws.Range('E13').Value = '''  -1.18%  '
$ws.Range('B14').Value = '''Chainlink'
$ws.Range('C14').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '''6.430'
$ws.Range('E14').Value = '''  +0.71%  '
$ws.Range('D15').Value = '''0.06871'
$ws.Range('E15').Value = '''  +0.43%  '
$ws.Range('D16').Value = '''1.004'
$ws.Range('E16').Value = '''  +0.23%  '
$ws.Range('D17').Value = '''78.78'
$ws.Range('E17').Value = '''  -2.75%  '
$ws.Range('D18').Value = '''0.000008684'
$ws.Range('E18').Value = '''  -0.84%  '
$ws.Range('D19').Value = '''1.004'
$ws.Range('E19').Value = '''  +0.20%  '
$ws.Range('D20').Value = '''14.88'
$ws.Range('E20').Value = '''  -1.78%  '
$ws.Range('D21').Value = '''26.437.54'
$ws.Range('E21').Value = '''  -2.70%  '
$ws.Range('D22').Value = '''5.092'
$ws.Range('E22').Value = '''  -0.41%  '
$ws.Range('D23').Value = '''11.07'
$ws.Range('E23').Value = '''  +1.96%  '
$ws.Range('D24').Value = '''2.015.17'
$ws.Range('E24').Value = '''  +0.47%  '
$ws.Range('D25').Value = '''152.25'
$ws.Range('E25').Value = '''  -0.95%  '
$ws.Range('D26').Value = '''1.812'
$ws.Range('E26').Value = '''  -8.15%  '
$ws.Range('D27').Value = '''17.98'
$ws.Range('E27').Value = '''  -0.79%  '
$ws.Range('D28').Value = '''5.071'
$ws.Range('E28').Value = '''  +0.02%  '
$ws.Range('D29').Value = '''113.80'
$ws.Range('E29').Value = '''  +0.39%  '
$ws.Range('D30').Value = '''1.784'
$ws.Range('E30').Value = '''  +4.99%  '
$ws.Range('E31').Value = '''  -0.02%  '
$ws.Range('D32').Value = '''0.7238'
$ws.Range('E32').Value = '''  -2.36%  '
$ws.Range('E33').Value = '''  +0.88%  '
$ws.Range('D34').Value = '''4.313'
$ws.Range('E34').Value = '''  -3.05%  '
$ws.Range('E35').Value = '''  +0.20%  '
$ws.Range('E36').Value = '''  -5.85%  '
$ws.Range('D37').Value = '''1.103'
$ws.Range('E37').Value = '''  +3.42%  '
$ws.Range('D38').Value = '''0.05134'
$ws.Range('E38').Value = '''  -0.92%  '
$ws.Range('D39').Value = '''0.01883'
$ws.Range('E39').Value = '''  -1.01%  '
$ws.Range('E40').Value = '''  -1.39%  '
$ws.Range('D41').Value = '''0.4910'
$ws.Range('E41').Value = '''  -1.05%  '
$ws.Range('D42').Value = '''2.602'
$ws.Range('E42').Value = '''  -3.68%  '
$ws.Range('D43').Value = '''6.344'
$ws.Range('E43').Value = '''  +1.07%  '
$ws.Range('D44').Value = '''7.951'
$ws.Range('E44').Value = '''  -2.35%  '
$ws.Range('D45').Value = '''104.72'
$ws.Range('E45').Value = '''  -0.34%  '
$ws.Range('D46').Value = '''1.004'
$ws.Range('E46').Value = '''  +0.28%  '
$ws.Range('D47').Value = '''10.18'
$ws.Range('E47').Value = '''  -0.51%  '
$ws.Range('D48').Value = '''1.639'
$ws.Range('E48').Value = '''  +3.07%  '
$ws.Range('D49').Value = '''0.06177'
$ws.Range('E49').Value = '''  -3.18%  '
$ws.Range('D50').Value = '''0.4426'
$ws.Range('E50').Value = '''  -3.02%  '
$ws.Range('D51').Value = '''1.714'
$ws.Range('E51').Value = '''  +0.58%  '
